# Auto-generated edit script: applies numeric corrections to market-price
# columns (H-N) across all 8 crafting-class Leve sheets, per scheduled
# Sheets runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 10
$ws.Range("H10").Value = 9995
$ws.Range("I10").Value = 9995
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 9995
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -9702
$ws.Range("N10").ClearContents()

# Row 33
$ws.Range("H33").Value = 818.93335
$ws.Range("I33").Value = 535
$ws.Range("J33").Value = 1599.75
$ws.Range("K33").Value = 535
$ws.Range("L33").Value = 1599.75
$ws.Range("M33").Value = -306
$ws.Range("N33").Value = -2057.75

# Row 40
$ws.Range("H40").Value = 5380.9585
$ws.Range("J40").Value = 6727
$ws.Range("L40").Value = 6727
$ws.Range("N40").Value = -7077

# Row 43
$ws.Range("H43").Value = 5540.857
$ws.Range("I43").Value = 1161.5714
$ws.Range("J43").Value = 7730.5
$ws.Range("K43").Value = 1161.5714
$ws.Range("L43").Value = 7730.5
$ws.Range("M43").Value = -1092.5714
$ws.Range("N43").Value = -7868.5

# Row 138
$ws.Range("H138").Value = 2163.262
$ws.Range("I138").Value = 1439.875
$ws.Range("J138").Value = 3127.7778
$ws.Range("K138").Value = 4319.625
$ws.Range("L138").Value = 9383.3334
$ws.Range("M138").Value = 820.375
$ws.Range("N138").Value = -19663.3334

$ws = $wb.Worksheets.Item("ARM")
# Row 16
$ws.Range("H16").Value = 1600
$ws.Range("I16").Value = 1400
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1400
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -1113
$ws.Range("N16").Value = -2574

# Row 45
$ws.Range("H45").Value = 831.625
$ws.Range("I45").Value = 807.9231
$ws.Range("J45").Value = 934.3333
$ws.Range("K45").Value = 807.9231
$ws.Range("L45").Value = 934.3333
$ws.Range("M45").Value = -430.9231
$ws.Range("N45").Value = -1688.3333

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1625.0526
$ws.Range("I20").Value = 1231.6364
$ws.Range("J20").Value = 2166
$ws.Range("K20").Value = 1231.6364
$ws.Range("L20").Value = 2166
$ws.Range("M20").Value = -984.6364000000001
$ws.Range("N20").Value = -2660

# Row 80
$ws.Range("H80").Value = 456.42856
$ws.Range("J80").Value = 297.8
$ws.Range("L80").Value = 297.8
$ws.Range("N80").Value = -2293.8

# Row 83
$ws.Range("H83").Value = 456.42856
$ws.Range("J83").Value = 297.8
$ws.Range("L83").Value = 1489
$ws.Range("N83").Value = -11473

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2069854.9
$ws.Range("I31").Value = 1007.7059
$ws.Range("J31").Value = 3372462.2
$ws.Range("K31").Value = 1007.7059
$ws.Range("L31").Value = 3372462.2
$ws.Range("M31").Value = -712.7059
$ws.Range("N31").Value = -3373052.2

# Row 34
$ws.Range("H34").Value = 2069854.9
$ws.Range("I34").Value = 1007.7059
$ws.Range("J34").Value = 3372462.2
$ws.Range("K34").Value = 1007.7059
$ws.Range("L34").Value = 3372462.2
$ws.Range("M34").Value = -805.7059
$ws.Range("N34").Value = -3372866.2

# Row 88
$ws.Range("H88").Value = 10000
$ws.Range("J88").Value = 10000
$ws.Range("L88").Value = 10000
$ws.Range("N88").Value = -10812

# Row 91
$ws.Range("H91").Value = 10000
$ws.Range("J91").Value = 10000
$ws.Range("L91").Value = 10000
$ws.Range("N91").Value = -12808

$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 5111.5
$ws.Range("I17").Value = 5111.5
$ws.Range("K17").Value = 15334.5
$ws.Range("M17").Value = -15165.5

# Row 95
$ws.Range("H95").Value = 8750
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 8750
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 26250
$ws.Range("M95").ClearContents()
$ws.Range("N95").Value = -30368

# Row 107
$ws.Range("H107").Value = 1096.3125
$ws.Range("I107").Value = 514.2
$ws.Range("J107").Value = 1360.909
$ws.Range("K107").Value = 1542.6
$ws.Range("L107").Value = 4082.727
$ws.Range("M107").Value = 377.3999999999999
$ws.Range("N107").Value = -7922.727000000001

# Row 131
$ws.Range("H131").Value = 1056.75
$ws.Range("I131").Value = 705.5789
$ws.Range("J131").Value = 1182.6415
$ws.Range("K131").Value = 2116.7367
$ws.Range("L131").Value = 3547.9245
$ws.Range("M131").Value = 2923.2633
$ws.Range("N131").Value = -13627.9245

$ws = $wb.Worksheets.Item("GSM")
# Row 9
$ws.Range("H9").Value = 907
$ws.Range("I9").Value = 907
$ws.Range("K9").Value = 907
$ws.Range("M9").Value = -737

# Row 70
$ws.Range("H70").Value = 6694.089
$ws.Range("I70").Value = 4380
$ws.Range("J70").Value = 7851.1333
$ws.Range("K70").Value = 4380
$ws.Range("L70").Value = 7851.1333
$ws.Range("M70").Value = -4110
$ws.Range("N70").Value = -8391.133300000001

# Row 73
$ws.Range("H73").Value = 6694.089
$ws.Range("I73").Value = 4380
$ws.Range("J73").Value = 7851.1333
$ws.Range("K73").Value = 4380
$ws.Range("L73").Value = 7851.1333
$ws.Range("M73").Value = -3444
$ws.Range("N73").Value = -9723.133300000001

# Row 126
$ws.Range("H126").Value = 2249.318
$ws.Range("I126").Value = 1687
$ws.Range("J126").Value = 3061.5557
$ws.Range("K126").Value = 5061
$ws.Range("L126").Value = 9184.667099999999
$ws.Range("M126").Value = -2591
$ws.Range("N126").Value = -14124.6671

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 28573866
$ws.Range("I2").Value = 1300.625
$ws.Range("J2").Value = 37039812
$ws.Range("K2").Value = 1300.625
$ws.Range("L2").Value = 37039812
$ws.Range("M2").Value = -1188.625
$ws.Range("N2").Value = -37040036

# Row 3
$ws.Range("H3").Value = 4166.6665
$ws.Range("J3").Value = 4500
$ws.Range("L3").Value = 4500
$ws.Range("N3").Value = -4724

# Row 15
$ws.Range("H15").Value = 4166.6665
$ws.Range("J15").Value = 4500
$ws.Range("L15").Value = 4500
$ws.Range("N15").Value = -4840

# Row 46
$ws.Range("H46").Value = 1206.1915
$ws.Range("I46").Value = 1058.1562
$ws.Range("J46").Value = 1522
$ws.Range("K46").Value = 1058.1562
$ws.Range("L46").Value = 1522
$ws.Range("M46").Value = -870.1561999999999
$ws.Range("N46").Value = -1898

# Row 116
$ws.Range("H116").Value = 22500
$ws.Range("J116").Value = 22500
$ws.Range("L116").Value = 22500
$ws.Range("N116").Value = -31678

$ws = $wb.Worksheets.Item("WVR")
# Row 42
$ws.Range("H42").Value = 8558.799999999999
$ws.Range("I42").Value = 6944
$ws.Range("J42").Value = 8962.5
$ws.Range("K42").Value = 6944
$ws.Range("L42").Value = 8962.5
$ws.Range("M42").Value = -6566
$ws.Range("N42").Value = -9718.5

# Row 80
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

# Row 83
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# Row 101
$ws.Range("H101").Value = 11821.2
$ws.Range("J101").Value = 11821.2
$ws.Range("L101").Value = 11821.2
$ws.Range("N101").Value = -18311.2

# Row 136
$ws.Range("H136").Value = 16749367
$ws.Range("I136").Value = 18056574
$ws.Range("K136").Value = 54169722
$ws.Range("M136").Value = -54167172
